$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.007.49"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.640.35"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "'215.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "'0.5089"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").Value = "'0.2581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").Value = "'0.06359"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").Value = "'19.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "'0.07746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'4.299"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "1.636.26"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "'0.5474"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "26.023.71"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'197.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'4.461"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'9.970"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'6.140"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'1.893"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").Value = "'142.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.45%  "
$ws.Range("D26").Value = "'0.1264"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.23%  "
$ws.Range("D27").Value = "'6.876"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").Value = "'15.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'0.04892"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("D31").Value = "'3.284"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").Value = "'3.214"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").Value = "'2.375"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'0.9187"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.5558"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.561"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "1.102.72"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").Value = "'0.01568"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "'5.615"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").Value = "'98.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  -4.62%  "
$ws.Range("D45").Value = "1.782.82"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'55.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "'0.05191"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("D50").Value = "'7.545"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("E51").Value = "  -0.19%  "
